# Scheduled runner update: refresh market-price-derived profit columns (H-N)
# across the per-job "Profits" worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 224.54546
$ws.Range("I11").Value = 224.54546
$ws.Range("K11").Value = 224.54546
$ws.Range("M11").Value = -84.54545999999999

$ws.Range("H61").Value = 3965.3333
$ws.Range("I61").Value = 5348
$ws.Range("K61").Value = 16044
$ws.Range("M61").Value = -15872

$ws.Range("H82").Value = 700
$ws.Range("I82").Value = 700
$ws.Range("K82").Value = 2100
$ws.Range("M82").Value = -1694

$ws.Range("H85").Value = 700
$ws.Range("I85").Value = 700
$ws.Range("K85").Value = 2100
$ws.Range("M85").Value = -696

$ws.Range("H87").Value = 133184.67
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 133184.67
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 133184.67
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -135680.67

$ws.Range("H90").Value = 133184.67
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 133184.67
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 399554.01
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -412034.01

$ws.Range("H99").Value = 313.4
$ws.Range("I99").Value = 313.4
$ws.Range("K99").Value = 940.1999999999999
$ws.Range("M99").Value = 557.8000000000001

$ws.Range("H101").Value = 849.6667
$ws.Range("I101").Value = 550
$ws.Range("K101").Value = 1650
$ws.Range("M101").Value = -28

$ws.Range("H104").Value = 799

$ws.Range("H112").Value = 1603.24
$ws.Range("I112").Value = 395
$ws.Range("J112").Value = 1768
$ws.Range("K112").Value = 1185
$ws.Range("L112").Value = 5304
$ws.Range("M112").Value = -77
$ws.Range("N112").Value = -7520

$ws.Range("H113").Value = 4803.4
$ws.Range("J113").Value = 6006
$ws.Range("L113").Value = 6006
$ws.Range("N113").Value = -12514

$ws.Range("H118").Value = 775.5714
$ws.Range("I118").Value = 738.3333
$ws.Range("K118").Value = 2214.9999
$ws.Range("M118").Value = -557.9998999999998

$ws.Range("H127").Value = 569
$ws.Range("I127").Value = 569
$ws.Range("K127").Value = 1707
$ws.Range("M127").Value = 3253

$ws.Range("H132").Value = 8458.117
$ws.Range("I132").Value = 8893
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 26679
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -24149
$ws.Range("N132").Value = -9560

$ws.Range("H135").Value = 500.72726
$ws.Range("I135").Value = 544.44446
$ws.Range("K135").Value = 4900.00014
$ws.Range("M135").Value = -2365.00014

$ws.Range("H138").Value = 2275.24
$ws.Range("I138").Value = 1101.8889
$ws.Range("J138").Value = 2935.25
$ws.Range("K138").Value = 3305.6667
$ws.Range("L138").Value = 8805.75
$ws.Range("M138").Value = 1834.3333
$ws.Range("N138").Value = -19085.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3595.3809
$ws.Range("J61").Value = 8666.333000000001
$ws.Range("L61").Value = 8666.333000000001
$ws.Range("N61").Value = -9090.333000000001

$ws.Range("H63").Value = 1750
$ws.Range("I63").Value = 1500
$ws.Range("K63").Value = 1500
$ws.Range("M63").Value = -814

$ws.Range("H66").Value = 1750
$ws.Range("I66").Value = 1500
$ws.Range("K66").Value = 7500
$ws.Range("M66").Value = -4068

$ws.Range("H97").Value = 796.7895
$ws.Range("I97").Value = 771.6111
$ws.Range("K97").Value = 771.6111
$ws.Range("M97").Value = -275.6111

$ws.Range("H102").Value = 2890.0588
$ws.Range("I102").Value = 2261.1667
$ws.Range("K102").Value = 2261.1667
$ws.Range("M102").Value = -639.1667000000002

$ws.Range("H122").Value = 2799.5454
$ws.Range("I122").Value = 2683.7856
$ws.Range("K122").Value = 8051.3568
$ws.Range("M122").Value = -5601.3568

$ws.Range("H132").Value = 2267.75
$ws.Range("I132").Value = 2267.75
$ws.Range("K132").Value = 6803.25
$ws.Range("M132").Value = -4273.25

$ws.Range("H136").Value = 3595.3809
$ws.Range("J136").Value = 8666.333000000001
$ws.Range("L136").Value = 25998.999
$ws.Range("N136").Value = -31098.999

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1179.4103
$ws.Range("I134").Value = 742.3939
$ws.Range("J134").Value = 3583
$ws.Range("K134").Value = 2227.1817
$ws.Range("L134").Value = 10749
$ws.Range("M134").Value = 307.8182999999999
$ws.Range("N134").Value = -15819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

$ws.Range("H122").Value = 628.1111
$ws.Range("I122").Value = 660.6667
$ws.Range("J122").Value = 563
$ws.Range("K122").Value = 1982.0001
$ws.Range("L122").Value = 1689
$ws.Range("M122").Value = 467.9999
$ws.Range("N122").Value = -6589

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 14066.667
$ws.Range("J23").Value = 20837.5
$ws.Range("L23").Value = 62512.5
$ws.Range("N23").Value = -62982.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5704.4287
$ws.Range("I122").Value = 4953.5
$ws.Range("J122").Value = 6004.8
$ws.Range("K122").Value = 14860.5
$ws.Range("L122").Value = 18014.4
$ws.Range("M122").Value = -12410.5
$ws.Range("N122").Value = -22914.4

$ws.Range("H123").Value = 76246.336
$ws.Range("J123").Value = 76246.336
$ws.Range("L123").Value = 76246.336
$ws.Range("N123").Value = -81146.336

$ws.Range("H126").Value = 7721
$ws.Range("I126").Value = 2397.7778
$ws.Range("J126").Value = 14565.143
$ws.Range("K126").Value = 7193.3334
$ws.Range("L126").Value = 43695.429
$ws.Range("M126").Value = -4723.3334
$ws.Range("N126").Value = -48635.429

$ws.Range("H132").Value = 3051.0908
$ws.Range("I132").Value = 2795
$ws.Range("K132").Value = 8385
$ws.Range("M132").Value = -5855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3244.75
$ws.Range("J7").Value = 3499.5
$ws.Range("L7").Value = 3499.5
$ws.Range("N7").Value = -3723.5

$ws.Range("H40").Value = 31089
$ws.Range("I40").Value = 39257.285
$ws.Range("K40").Value = 39257.285
$ws.Range("M40").Value = -39121.285

$ws.Range("H126").Value = 3244.75
$ws.Range("J126").Value = 3499.5
$ws.Range("L126").Value = 10498.5
$ws.Range("N126").Value = -15438.5

$ws.Range("H132").Value = 4219.1113
$ws.Range("I132").Value = 3912
$ws.Range("J132").Value = 4833.3335
$ws.Range("K132").Value = 11736
$ws.Range("L132").Value = 14500.0005
$ws.Range("M132").Value = -9206
$ws.Range("N132").Value = -19560.0005

$ws.Range("H136").Value = 5912
$ws.Range("I136").Value = 3824.25
$ws.Range("K136").Value = 11472.75
$ws.Range("M136").Value = -8922.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2542.5
$ws.Range("I132").Value = 2346.8572
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 7040.571599999999
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -4510.571599999999
$ws.Range("N132").Value = -14057

